$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B column: repurpose "Time" header as "Time (in minutes)" and add minute counts ---
$ws.Range("B1").Value = "Time (in minutes)"
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 30
$ws.Range("B4").Value = 205

# --- New row 5: 06.11.2023 / 150 min / battle logic / 13:00-15:30 ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "06.11.2023"
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("B5").Value = 150
$ws.Range("C5").Value = "Worked on battle logic"

$ws.Range("D5").Value = "13:00-15:30"
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# --- New row 6: 08.12.2023 / 270 min / HTTP server / 15:30-20:00 ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08.12.2023"
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B6").Value = 270
$ws.Range("C6").Value = "Got started on HTTP Server"

$ws.Range("D6").Value = "15:30-20:00"
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)

# --- New row 7: stray note, only in column D ---
$ws.Range("D7").Value = "13:30-zeit"
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# --- New column F: running hour counter ---
$ws.Range("F1").Value = "Temp Hour counter"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Formula = "=SUM(B2:B250)/60"

# --- Column widths ---
$ws.Columns("B:B").ColumnWidth = 14
$ws.Columns("F:F").ColumnWidth = 16

# --- Clear clipboard marquee / selection state ---
$excel.CutCopyMode = $false
[void]$ws.Range("D8").Select()
